$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Ticker -> KHC (Quarter/Year stay the same: q1 / 2019)
$ws.Range("A2").Value = "KHC"

# Update row 3: Ticker -> MCHP, clear Quarter/Year
$ws.Range("A3").Value = "MCHP"
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Remove old rows 4-12 entirely so the used range shrinks back to A1:C3
$ws.Range("A4:C12").Delete(-4162)

$ws.Range("C8").Select() | Out-Null
